# Apply the process-flow renaming changes described in the diff.
# Only the B and C column text values for the relevant rows change;
# styles/formatting remain untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value  = "SelfCareInPerizia"
$ws.Range("B4").Value  = "SelfCareInChiusura"

$ws.Range("B25").Value = "AppuntamentoModificato"
$ws.Range("C25").Value = "EsecuzioneDesk"

$ws.Range("B26").Value = "AppuntamentoAnnullato"

$ws.Range("B27").Value = "DatiObbligatoriMancanti"

$ws.Range("B28").Value = "DeskRifiutata"
$ws.Range("C28").Value = "RiassegnazioneDesk"

$ws.Range("B29").Value = "DeskCompletata"
$ws.Range("C29").Value = "AvvioPeriziaPostDesk"
